$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
$ws.Range("A8").Value = "Volume 32   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/13/2025  Through  1/19/2025"

# --- Numeric value updates across the weekly crime-stats table ---
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 100
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 11
$ws.Range("K16").Value = -45.454545454545
$ws.Range("L16").Value = -33.333333333333
$ws.Range("M16").Value = -40
$ws.Range("N16").Value = -93.103448275862
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 13
$ws.Range("H17").Value = 62.5
$ws.Range("I17").Value = 7
$ws.Range("J17").Value = 4
$ws.Range("K17").Value = 75
$ws.Range("L17").Value = -12.5
$ws.Range("M17").Value = 75
$ws.Range("N17").Value = -22.222222222222
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 84.615384615384
$ws.Range("I18").Value = 15
$ws.Range("J18").Value = 10
$ws.Range("K18").Value = 50
$ws.Range("L18").Value = 15.384615384615
$ws.Range("M18").Value = -31.818181818181
$ws.Range("N18").Value = -88.721804511278
$ws.Range("C19").Value = 30
$ws.Range("D19").Value = 34
$ws.Range("E19").Value = -11.764705882352
$ws.Range("F19").Value = 118
$ws.Range("G19").Value = 131
$ws.Range("H19").Value = -9.923664122137
$ws.Range("I19").Value = 77
$ws.Range("J19").Value = 83
$ws.Range("K19").Value = -7.22891566265
$ws.Range("L19").Value = 14.925373134328
$ws.Range("M19").Value = 20.3125
$ws.Range("N19").Value = -58.823529411764
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 2
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -66.666666666666
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -98.963730569948
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -11.363636363636
$ws.Range("F21").Value = 171
$ws.Range("G21").Value = 171
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 109
$ws.Range("J21").Value = 111
$ws.Range("K21").Value = -1.801801801801
$ws.Range("L21").Value = 4.807692307692
$ws.Range("M21").Value = 5.825242718446
$ws.Range("N21").Value = -82.218597063621
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = -50
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 0
$ws.Range("M23").Value = -50
$ws.Range("C24").Value = 65
$ws.Range("D24").Value = 77
$ws.Range("E24").Value = -15.584415584415
$ws.Range("F24").Value = 258
$ws.Range("G24").Value = 223
$ws.Range("H24").Value = 15.695067264574
$ws.Range("I24").Value = 176
$ws.Range("J24").Value = 159
$ws.Range("K24").Value = 10.691823899371
$ws.Range("L24").Value = 15.032679738562
$ws.Range("M24").Value = 112.048192771084
$ws.Range("C25").Value = 56
$ws.Range("D25").Value = 71
$ws.Range("E25").Value = -21.12676056338
$ws.Range("F25").Value = 220
$ws.Range("G25").Value = 215
$ws.Range("H25").Value = 2.325581395348
$ws.Range("I25").Value = 143
$ws.Range("J25").Value = 149
$ws.Range("K25").Value = -4.026845637583
$ws.Range("L25").Value = 5.925925925925
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = -18.181818181818
$ws.Range("I26").Value = 13
$ws.Range("J26").Value = 15
$ws.Range("K26").Value = -13.333333333333
$ws.Range("L26").Value = -27.777777777777
$ws.Range("M26").Value = -38.095238095238
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("L27").Value = 100
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -12.5
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 4
$ws.Range("K28").Value = 75
$ws.Range("G31").Value = 3
$ws.Range("J43").Value = 255
$ws.Range("K43").Value = -63.623395149786
$ws.Range("L43").Value = -75.194552529182
$ws.Range("M43").Value = -91.553494534614
$ws.Range("N43").Value = -93.098782138024
$ws.Range("J46").Value = 2449
$ws.Range("K46").Value = -34.255033557047
$ws.Range("L46").Value = -49.598682856554
$ws.Range("M46").Value = -79.550768203072
$ws.Range("N46").Value = -83.329929889047

# --- Cells that become the text "0" (insufficient-data marker) ---
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"

# --- Cells that become the text "***.*" (undefined pct-change marker) ---
$ws.Range("E15").Value = "***.*"
$ws.Range("E27").Value = "***.*"
$ws.Range("E31").Value = "***.*"
